$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''27.376.64'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.66%  '
$ws.Range("D3").Value = '''1.840.31'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.34%  '
$ws.Range("D4").Value = '''1.015'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +1.35%  '
$ws.Range("D5").Value = '''315.03'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.96%  '
$ws.Range("E6").Value = '  +1.17%  '
$ws.Range("E7").Value = '  +1.41%  '
$ws.Range("E8").Value = '  +0.30%  '
$ws.Range("D9").Value = '''0.07461'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.31%  '
$ws.Range("D10").Value = '''0.8856'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.72%  '
$ws.Range("D11").Value = '''20.51'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.52%  '
$ws.Range("D12").Value = '''1.848.86'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.24%  '
$ws.Range("D13").Value = '''0.07373'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +4.20%  '
$ws.Range("D14").Value = '''5.485'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.92%  '
$ws.Range("D15").Value = '''93.29'
$ws.Range("D15").Style = "Normal"
$ws.Range("E16").Value = '  +1.02%  '
$ws.Range("E17").Value = '  +1.24%  '
$ws.Range("D18").Value = '''0.000008853'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.66%  '
$ws.Range("D19").Value = '''1.013'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.22%  '
$ws.Range("D20").Value = '''14.84'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.89%  '
$ws.Range("D21").Value = '''27.396.30'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.64%  '
$ws.Range("D22").Value = '''5.353'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.53%  '
$ws.Range("D23").Value = '''10.73'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.19%  '
$ws.Range("D24").Value = '''2.067.63'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.28%  '
$ws.Range("E25").Value = '  +0.68%  '
$ws.Range("D26").Value = '''152.37'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.12%  '
$ws.Range("D27").Value = '''18.65'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.39%  '
$ws.Range("D28").Value = '''2.168'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.22%  '
$ws.Range("D29").Value = '''5.264'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.37%  '
$ws.Range("D30").Value = '''118.12'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.89%  '
$ws.Range("E31").Value = '  +0.21%  '
$ws.Range("D32").Value = '''0.7611'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.15%  '
$ws.Range("D33").Value = '''1.178'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.37%  '
$ws.Range("D34").Value = '''4.562'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.17%  '
$ws.Range("D35").Value = '''2.943'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.10%  '
$ws.Range("E36").Value = '  +1.26%  '
$ws.Range("D37").Value = '''1.107'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.78%  '
$ws.Range("D38").Value = '''0.05384'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.74%  '
$ws.Range("D39").Value = '''0.01963'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.03%  '
$ws.Range("D40").Value = '''3.001'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.34%  '
$ws.Range("D41").Value = '''7.312'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.57%  '
$ws.Range("D42").Value = '''0.5357'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.71%  '
$ws.Range("D43").Value = '''2.384'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.99%  '
$ws.Range("D44").Value = '''0.1668'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.61%  '
$ws.Range("D45").Value = '''8.562'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.47%  '
$ws.Range("D46").Value = '''0.4990'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.33%  '
$ws.Range("D47").Value = '''10.58'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.58%  '
$ws.Range("E48").Value = '  +1.35%  '
$ws.Range("D49").Value = '''105.24'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.35%  '
$ws.Range("D50").Value = '''1.681'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.63%  '
$ws.Range("D51").Value = '''0.06324'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.36%  '
